# Release MHD 4.2.2 close #419
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 4.2.1 -> 4.2.2
$ws.Range("B3").Value = "4.2.2"

# Date: refresh publication date
$ws.Range("B8").Value = "2024-05-18T12:39:23-05:00"

# Contact rows: replace placeholder "No display for ContactDetail" text
# with the actual rendered ContactDetail values (url, email, name+email)
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"
